$wb = $excel.ActiveWorkbook

# --- Sheet "Question Validation Succeed" ---
$wsSucceed = $wb.Worksheets.Item("Question Validation Succeed")
$wsSucceed.Activate()
$wsSucceed.Range("O3").Value = "1+1"
$wsSucceed.Range("O4").Value = "2+2"
$wsSucceed.Range("P5").Value = '{ "column": "xyz", "writeToPatient": { "fieldName": "abc", "fieldType": "FreeText" } }'
$wsSucceed.Range("P5").Select() | Out-Null

# --- Sheet "Question Validation Fail" ---
$wsFail = $wb.Worksheets.Item("Question Validation Fail")
$wsFail.Activate()
$wsFail.Range("O6").Value = "1+1"
$wsFail.Range("O7").Value = "2+2"
$wsFail.Range("O8").Value = "1+1"
$wsFail.Range("O9").Value = "2+2"
$wsFail.Range("O10").Value = "1+1"
$wsFail.Range("O11").Value = "2+2"
$wsFail.Range("O12").Value = "1+1"
$wsFail.Range("O13").Value = "2+2"
$wsFail.Range("O6").Select() | Out-Null

$wsSucceed.Activate() | Out-Null
